$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph. It needs the same leading-empty-run / trailing
#    text-run shape used by the normal body paragraphs elsewhere in the
#    doc, so we seed it by pasting a copy of an existing body paragraph
#    (which already has that <w:r/><w:r>text</w:r> shape) and then
#    overwrite its text/formatting in place.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$bodyPara = $d.Paragraphs.Item(4)
$bodyPara.Range.Copy()

$metaRange = $metaPara.Range
$metaRange.Collapse(1)
$metaRange.Paste()

$metaPara = $d.Paragraphs.Item(2)
$full = $metaPara.Range

# Bold + overwrite the first 16 characters with the "Meta description" label.
$labelRange = $d.Range($full.Start, $full.Start + 16)
$labelRange.Bold = 1
$labelRange.Text = "Meta description"

# Overwrite the remainder (non-bold) with the rest of the sentence.
$full2 = $metaPara.Range
$restRange = $d.Range($full.Start + 16, $full2.End - 1)
$restRange.Text = ": Experience terror with Castle of Terror, a horror-themed online slot game by Big Time Gaming. Play for free and read our review here to learn more."

# ---------------------------------------------------------------------
# 2) Near the end of the document: remove the duplicated bold title
#    paragraph ("Play Castle of Terror for Free - ...") entirely, and
#    replace the text of the remaining italic paragraph with the new
#    "Prompt: ..." text (keeping its italic formatting intact).
# ---------------------------------------------------------------------
$oldTitleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Castle of Terror for Free - Review of Big Time Gaming's Horror-themed Slot`r") {
        $oldTitleIndex = $i
    }
}

if ($oldTitleIndex -ge 1) {
    $d.Paragraphs.Item($oldTitleIndex).Range.Delete()
}

$promptIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Experience terror with Castle of Terror, a horror-themed online slot game by Big Time Gaming. Play for free and read our review here to learn more.`r") {
        $promptIndex = $i
    }
}

$promptPara = $d.Paragraphs.Item($promptIndex)
$fullRange = $promptPara.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$textRange.Text = "Prompt: Create a feature graphic for `"Castle of Terror`" Design a cartoon-style feature graphic that showcases a happy Maya warrior with glasses to fit the theme of `"Castle of Terror.`" The warrior should be holding a sword, and there should be spooky elements in the background such as a haunted castle and a full moon. The overall color scheme should be dark with pops of bright colors to add contrast and make the image pop. Add the title of the game `"Castle of Terror`" in a spooky and eye-catching font. The image should be in a square format so that it can be easily used on social media platforms as well."

Write-Output "done"
